$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("C","D","E","F","G","I","J","L","M","N")

$data = @(
    @(5.213037744005318,8.560862889749792,7.667890652881365,79.98064016098319,3.91386196359629,60.58479843161717,9.944455771638804,8.417923520109101,44.57834913285934,17.66954877468222),
    @(5.066173086643301,8.454458234354824,7.555733644705913,79.81284328313653,3.923217282932987,60.4952462285472,9.985375360664314,8.419756092376508,44.002219550803,17.56914069563858),
    @(4.97621289360604,8.392082455088016,7.485420822034916,79.73835262761789,3.929231478432391,60.461302787966,10.01157862991889,8.421231320146799,43.65384367868774,17.50992574574658),
    @(4.939665681912286,8.367428652280035,7.45641672784768,79.71512426051143,3.931750674135258,60.45272183896672,10.02252958851129,8.421920415839805,43.51336702810109,17.48642872259476),
    @(4.933605565157587,8.363381681313662,7.451579806257133,79.71169610122708,3.932173127270113,60.45161288474129,10.02436452024695,8.422040147824534,43.49013495209302,17.48256592445081),
    @(4.975719472136615,8.391746841937394,7.485031066556576,79.73801057379832,3.92926517576821,60.46116585612757,10.01172521108214,8.421240257623765,43.65194294931965,17.50960626422132),
    @(5.162388353455708,8.52357417043892,7.62952714459629,79.91683377648845,3.917031926286954,60.54953098087602,9.958342054191679,8.418482689612919,44.37866917741282,17.6344352703974),
    @(5.527721191824797,8.804526386598988,7.90084645536739,80.49577707977055,3.89516318360377,60.89134347764574,9.862135215444731,8.415857052066182,45.84019867202672,17.89760249161515),
    @(5.79264775274407,9.023147578842242,8.092152797086348,81.06231868796505,3.88035842458628,61.24705110196672,9.796506290813637,8.415631577945899,46.92787532528829,18.10084579383131),
    @(5.911840072548525,9.124919103993848,8.177302869389047,81.35098594520598,3.873890720170163,61.43186938734493,9.76772244013911,8.415900732718553,47.42395665781991,18.19516708774717),
    @(5.956742094836668,9.163761912186949,8.209266418932277,81.46475585512485,3.871479428592113,61.50517843472427,9.756974762520571,8.416056250784605,47.61185220677175,18.23112754902795),
    @(5.947082671163971,9.15538336115501,8.202395161943359,81.44005515484623,3.871997066597164,61.48924208221796,9.75928273089324,8.416020371202269,47.57138602126886,18.22337241322074),
    @(5.915539101880643,9.128108798627038,8.179938221739993,81.36025650212908,3.873691585415896,61.43783391769367,9.766835183966005,8.415912452324385,47.4394149496561,18.19812085035053),
    @(5.896186131621022,9.111441074472848,8.166145791557875,81.31195818264514,3.874734445798036,61.40677794145536,9.771481036618738,8.415853332870688,47.35857985424445,18.18268444028608),
    @(5.784827833477464,9.016540736138795,8.08654931933903,81.04407715249866,3.880786434327346,61.23543704822288,9.798408776182967,8.415621480018988,46.89546872310126,18.09471705751103),
    @(5.716142052420578,8.958895698882708,8.037230804077979,80.88767463959765,3.884567169083113,61.13622886925866,9.815201074785385,8.415574551675316,46.611591315277,18.0412117407443),
    @(5.676513504962458,8.925959917389987,8.008688134531289,80.80062871454466,3.886766915453609,61.08133176927998,9.824960489792232,8.415582539285897,46.44844263142988,18.01061394959887),
    @(5.663076332926509,8.914847126912937,7.998994205148858,80.77165634231335,3.887516049892158,61.06311584724333,9.828282251351538,8.415591247329075,46.3932301230456,18.00028523481414),
    @(5.723466753920666,8.965009542670851,8.042499138600478,80.90402242200977,3.884162102348306,61.14656557930707,9.813403073665597,8.415575926025701,46.64179815084304,18.04688933035947),
    @(5.924810874481561,9.13611197950873,8.186542079611332,81.38357428768734,3.873192840042194,61.45284348445469,9.764612731452289,8.415942694995984,47.47817809884412,18.20553145172078),
    @(6.055022555280361,9.249696425635932,8.279040253246391,81.72297163481548,3.866244386500316,61.67238225913469,9.733611327002697,8.416494863523383,48.02497106591754,18.31061808248362),
    @(5.985665466383084,9.188923041159919,8.22982601557889,81.53945073777876,3.869932894890364,61.55343463891223,9.750076928599949,8.41617152158941,47.73316883248236,18.25441116980249),
    @(5.72015568967665,8.962244832765387,8.040117911773907,80.89662264475676,3.884345151524174,61.1418856913145,9.814215621847865,8.415575195777038,46.62814145430023,18.04432198208313),
    @(5.429293584038683,8.726266412170878,7.82880963919682,80.31447897685925,3.900855366451361,60.78062226376463,9.88726581933253,8.416268759586535,45.44177673871702,17.8245745234673)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $rowValues = $data[$i]
    for ($j = 0; $j -lt $cols.Count; $j++) {
        $ws.Range("$($cols[$j])$row").Value = $rowValues[$j]
    }
}
